$wb = $excel.ActiveWorkbook

# Map sheet name -> hashtable of cell => new value
$updates = @{
    "2025" = @{
        "B2" = 1037.265132737054
        "E2" = 28926.05393052954
        "G2" = 8095.925712661834
        "I2" = 16171.06685703679
        "L2" = 48492.22142001599
        "M2" = 10595.37713982
        "N2" = 7068.711122921395
        "O2" = 6993.065970389833
    }
    "2030" = @{
        "A2" = 0
        "B2" = 4157.588990853394
        "E2" = 45991.90904307188
        "G2" = 8095.925712661834
        "I2" = 37079.12819938764
        "L2" = 54844.03303316472
        "M2" = 17449.04999683176
        "N2" = 9020.386661498747
        "O2" = 9721.982264164202
    }
    "2035" = @{
        "A2" = 2754.31755456332
        "B2" = 6368.910634126893
        "E2" = 57457.45307013817
        "G2" = 8095.925712661834
        "I2" = 52465.73681402855
        "L2" = 54844.03303316472
        "M2" = 21912.87293902603
        "N2" = 13027.78294739439
        "O2" = 12857.64537285333
    }
    "2040" = @{
        "A2" = 2754.31755456332
        "B2" = 6368.910634126893
        "E2" = 57457.45307013817
        "G2" = 8095.925712661834
        "I2" = 52465.73681402855
        "L2" = 54844.03303316472
        "M2" = 21912.87293902603
        "N2" = 13145.17518818184
        "O2" = 12857.64537285333
    }
    "2045" = @{
        "A2" = 5713.151062849596
        "B2" = 6368.910634126893
        "E2" = 57457.45307013817
        "G2" = 8095.925712661834
        "I2" = 52465.73681402855
        "L2" = 54844.03303316472
        "M2" = 21912.87293902603
        "N2" = 13593.64654746441
        "O2" = 14932.85150547986
    }
    "2050" = @{
        "A2" = 5713.151062849596
        "B2" = 6368.910634126893
        "E2" = 57457.45307013817
        "G2" = 8095.925712661834
        "I2" = 52465.73681402855
        "L2" = 54844.03303316472
        "M2" = 21912.87293902603
        "N2" = 13593.64654746441
        "O2" = 14932.85150547986
    }
}

foreach ($sheetName in $updates.Keys) {
    $sheetNameStr = [string]$sheetName
    $ws = $wb.Worksheets.Item($sheetNameStr)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range([string]$cellRef).Value = $cellUpdates[$cellRef]
    }
}
